# Update pl_mw.xlsx results for the 380 kV case.
# Columns B, C, E, F, G, J, M, O in rows 2-25 are updated with new
# computed values (columns A, D, H, I, K, L, N are unchanged zeros/indices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: @(rowNumber, B, C, E, F, G, J, M, O)
$data = @(
    @(2, 0.8028798507536976, 0.1892926474504861, 0.5136124957476795, 1.949043443498852, 0.002407936862410621, 0.03392508486171053, 0.573804906084554, 1.88418337364368),
    @(3, 0.7057052289942476, 0.1677811970580194, 0.5088575932777246, 1.944057408184634, 0.002410781839193269, 0.03437957001682435, 0.5351438526236123, 1.906550414271592),
    @(4, 0.6459099745758863, 0.1545129747835574, 0.5061439439186799, 1.942305317627927, 0.00241262027346376, 0.0346803107751299, 0.5115171036108208, 1.922159333167841),
    @(5, 0.6215115727715954, 0.1490912833723428, 0.5050899413793815, 1.941920121427032, 0.002413392555926098, 0.03480831920645855, 0.5019174584617687, 1.928990246737342),
    @(6, 0.617458380412927, 0.1481901326066861, 0.5049180573982213, 1.941876005079678, 0.002413522190349779, 0.03482990436524247, 0.5003251765225087, 1.930152868283358),
    @(7, 0.645581054461843, 0.154439915370773, 0.5061295193258815, 1.94229879206415, 0.002412630595093433, 0.03468201505666002, 0.5113875236429593, 1.92224955555848),
    @(8, 0.7694016642230395, 0.1818881647349428, 0.5119303089906992, 1.947052221501309, 0.002408898839725395, 0.03407729132948134, 0.5604517382455469, 1.89150540523643),
    @(9, 1.01114229279176, 0.2352253623486433, 0.524937647235781, 1.9667867988654, 0.002402304527842885, 0.03306344278396978, 0.6575345969959585, 1.846163049023119),
    @(10, 1.188057138679198, 0.2741023260790882, 0.5354886502608451, 1.98767219420526, 0.002397896327370333, 0.03242337738478618, 0.7293782204537109, 1.822053892362703),
    @(11, 1.268383135648833, 0.2917189625557057, 0.5405045670952191, 1.99856883086953, 0.00239598477268993, 0.03215494660122786, 0.7621719206449313, 1.813104335041373),
    @(12, 1.298777482549156, 0.2983797830935089, 0.5424350263513773, 2.002896412550356, 0.002395274326181346, 0.03205656878652796, 0.774605766017018, 1.810007127979048),
    @(13, 1.292232575766604, 0.296945715992365, 0.5420178876749233, 2.001955430422285, 0.002395426737799764, 0.0320776107017835, 0.7719272294131798, 1.81066116320369),
    @(14, 1.27088417250269, 0.2922671589707306, 0.5406627653303389, 1.998920826864321, 0.002395926055161156, 0.03214678743310273, 0.7631945509897378, 1.81284366856346),
    @(15, 1.257804578048081, 0.2894000670969774, 0.5398367546583032, 1.997088271208156, 0.002396233647090622, 0.03218958621586765, 0.7578475486205889, 1.814218565083792),
    @(16, 1.182804446600755, 0.272949626076894, 0.5351651946905136, 1.986988206953598, 0.002398023133053048, 0.03244137763097577, 0.7272372813997805, 1.822679510588358),
    @(17, 1.136754062699765, 0.2628399769207022, 0.5323546847746172, 1.981150026588537, 0.002399144892346165, 0.03260166821590893, 0.7084871362073102, 1.828388003574929),
    @(18, 1.110252732586105, 0.2570187287611532, 0.530758507723057, 1.977923402632243, 0.002399798927587656, 0.03269600324994038, 0.6977130881347335, 1.831861206786499),
    @(19, 1.101277409815452, 0.2550466575984842, 0.5302215669683008, 1.97685346467027, 0.002400021890828008, 0.03272831100210283, 0.6940670053507461, 1.833069727144419),
    @(20, 1.14165770110219, 0.2639168358819575, 0.5326517621697207, 1.981757913345604, 0.002399024565853565, 0.0325843835042221, 0.7104820338388862, 1.827760666289123),
    @(21, 1.277155360411825, 0.2936416445586758, 0.541059955582412, 1.999806696892733, 0.00239577902983915, 0.03212637976192667, 0.7657591311773047, 1.812194681547936),
    @(22, 1.36557413504039, 0.3130088164789981, 0.5467361005835727, 2.012775985482207, 0.002393736064394901, 0.03184611507586332, 0.801976521579931, 1.803722978724949),
    @(23, 1.31839631076565, 0.3026777667369913, 0.5436900988525863, 2.005746486330693, 0.002394819301262959, 0.03199395245387926, 0.7826384834946509, 1.808088252525692),
    @(24, 1.139440848167055, 0.2634300161944907, 0.5325173923526663, 1.981482683386048, 0.002399078937053293, 0.03259219112588063, 0.7095801229554866, 1.828043689655402),
    @(25, 0.9458638229712619, 0.2208498121540288, 0.5212441459287334, 1.960329322434745, 0.00240401146087202, 0.03331930693701857, 0.6311794773176302, 1.856819688444347),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]   # B
    $ws.Cells.Item($r, 3).Value  = $row[2]   # C
    $ws.Cells.Item($r, 5).Value  = $row[3]   # E
    $ws.Cells.Item($r, 6).Value  = $row[4]   # F
    $ws.Cells.Item($r, 7).Value  = $row[5]   # G
    $ws.Cells.Item($r, 10).Value = $row[6]   # J
    $ws.Cells.Item($r, 13).Value = $row[7]   # M
    $ws.Cells.Item($r, 15).Value = $row[8]   # O
}
